# Turn the single-sheet "Antibodies" submission workbook into a three-sheet
# workbook: Instructions, Antibodies, Terminology.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The existing sheet becomes "Antibodies" (keep its data, fix the header
#    typo, add dropdown validations sourced from the new Terminology sheet).
# ---------------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item(1)
$antibodies.Name = "Antibodies"

# "Isoform" -> "Isotype"
$antibodies.Range("C1").Value = "Isotype"

# ---------------------------------------------------------------------------
# 2. Add the "Instructions" sheet before "Antibodies" and fill it in.
# ---------------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item("Antibodies")
$instructions = $wb.Worksheets.Add($antibodies)
$instructions.Name = "Instructions"

$instructions.Cells.Item(1, 1).Value = "CoVIC-DB Antibodies Submission"
$instructions.Cells.Item(3, 1).Value = "Add your antibodies to the 'Antibodies' sheet. Do not edit the other sheets."
$instructions.Cells.Item(5, 1).Value = "Columns:"
$instructions.Cells.Item(6, 1).Value = "- Antibody name: Your institutions preferred name for the antibody."
$instructions.Cells.Item(7, 1).Value = "- Host: The name of the host species that is the source of the antibody."
$instructions.Cells.Item(8, 1).Value = "- Isotype: The name of the isotype of the antibody's heavy chain."

[void]$instructions.Range("A1").Select()
$instructions.PageSetup.LeftMargin = 54
$instructions.PageSetup.RightMargin = 54
$instructions.PageSetup.TopMargin = 72
$instructions.PageSetup.BottomMargin = 72
$instructions.PageSetup.HeaderMargin = 36
$instructions.PageSetup.FooterMargin = 36

$instructions.Protect()

# ---------------------------------------------------------------------------
# 3. Add the "Terminology" sheet after "Antibodies" and fill it in with the
#    reference lists of hosts / isotypes.
# ---------------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item("Antibodies")
$terminology = $wb.Worksheets.Add($null, $antibodies)
$terminology.Name = "Terminology"

$terminology.Cells.Item(1, 1).Value = "Host"
$terminology.Cells.Item(1, 2).Value = "Isotype"

$hosts = @("Homo sapiens", "Mus musculus")
for ($i = 0; $i -lt $hosts.Length; $i++) {
    $terminology.Cells.Item($i + 2, 1).Value = $hosts[$i]
}

$isotypes = @("IgA", "IgA1", "IgA2", "IgD", "IgE", "IgG", "IgG1", "IgG2", "IgG2a", "IgG2b", "IgG2c", "IgG3", "IgG4", "IgM", "sIgA")
for ($i = 0; $i -lt $isotypes.Length; $i++) {
    $terminology.Cells.Item($i + 2, 2).Value = $isotypes[$i]
}

$terminology.Range("A1:B1").Font.Bold = $true
# 14.14 chars of "ColumnWidth" rounds to a stored width of exactly 15.
$terminology.Columns.Item(1).ColumnWidth = 14.14
$terminology.Columns.Item(2).ColumnWidth = 14.14

[void]$terminology.Range("A2").Select()
$terminology.Application.ActiveWindow.FreezePanes = $true
[void]$terminology.Range("A1").Select()

$terminology.PageSetup.LeftMargin = 54
$terminology.PageSetup.RightMargin = 54
$terminology.PageSetup.TopMargin = 72
$terminology.PageSetup.BottomMargin = 72
$terminology.PageSetup.HeaderMargin = 36
$terminology.PageSetup.FooterMargin = 36

$terminology.Protect()

# ---------------------------------------------------------------------------
# 4. Back on "Antibodies": freeze-pane selection moves to A2, and the two
#    dropdown data validations are wired up to the Terminology sheet.
# ---------------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item("Antibodies")

$antibodies.Range("B2:B100").Validation.Add(3, 1, 1, "=Terminology!A2:A3")
$antibodies.Range("C2:C100").Validation.Add(3, 1, 1, "=Terminology!B2:B16")

[void]$antibodies.Range("A2").Select()
$antibodies.Application.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 5. "Antibodies" is the active tab (index 1, zero-based) when the workbook
#    is opened.
# ---------------------------------------------------------------------------
$antibodies = $wb.Worksheets.Item("Antibodies")
[void]$antibodies.Select()
